$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.385.82"
$ws.Range("E2").Value = "  -2.69%  "
$ws.Range("D3").Value = "3.781.84"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'593.44"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "'165.64"
$ws.Range("E6").Value = "  -3.01%  "
$ws.Range("D7").Value = "3.783.01"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("E11").Value = "  -2.37%  "
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("E13").Value = "  -3.93%  "
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("D15").Value = "4.417.64"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "3.771.52"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "67.388.76"
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("D18").Value = "'17.99"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "'10.19"
$ws.Range("E21").Value = "  -7.82%  "
$ws.Range("D22").Value = "'457.21"
$ws.Range("E22").Value = "  -3.08%  "
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").Value = "'83.43"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("D26").Value = "'11.83"
$ws.Range("E26").Value = "  -3.43%  "
$ws.Range("D27").Value = "'2.12"
$ws.Range("E27").Value = "  -5.19%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'9.94"
$ws.Range("E29").Value = "  -3.45%  "
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("D31").Value = "'29.75"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("D32").Value = "'2.18"
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("D33").Value = "'7.17"
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("D34").Value = "'9.15"
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("D36").Value = "3.735.50"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "'0.0995"
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("E39").Value = "  -7.54%  "
$ws.Range("D40").Value = "'0.990"
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("D41").Value = "'5.71"
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'43.82"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "'0.297"
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("D47").Value = "'8.33"
$ws.Range("E47").Value = "  -3.86%  "
$ws.Range("D48").Value = "'147.68"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").Value = "'392.14"
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("E50").Value = "  -8.12%  "
$ws.Range("D51").Value = "2.749.77"
$ws.Range("E51").Value = "  +1.86%  "
